# Update settings for system modelling
#
# 1. Switch workbook calculation to manual.
# 2. On the "rel_connection__node__node" sheet, change the
#    fix_ratio_out_in_connection_flow "value" column (F2:F243) from 1 to 0.99.

$wb = $excel.ActiveWorkbook
$excel.Calculation = -4135   # xlCalculationManual

$ws = $wb.Worksheets.Item("rel_connection__node__node")

$rng = $ws.Range("F2:F243")
$rng.Value = 0.99
